# Publishing after approval to publish (end-or-change-of-employment):
# Convert the page's one-and-only Heading1 into the page's title line
# (styled FirstParagraph, "#"-prefixed) and turn every Heading2 section
# header into a "##"-prefixed BodyText/FirstParagraph line, while the old
# FirstParagraph "intro" line that used to follow each heading becomes an
# ordinary BodyText paragraph. This mirrors the markdown-ish "#"/"##"
# heading convention used elsewhere on the intranet.

$d = $word.ActiveDocument

function Get-ParaByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $needle) {
            return $p
        }
    }
    return $null
}

# --- Section: page title -------------------------------------------------
$pTitle = Get-ParaByText "End or change of employment"
$pTitle.Range.Text = "#End or change of employment"
$pTitle.Style = "FirstParagraph"

$pIntro = Get-ParaByText "Managers must ensure that all employees, contractors and third-party users return all assets within their possession and that all access rights (including building passes, access to buildings, IT systems, applications and directories) are removed at the point of termination or change of employment."
$pIntro.Style = "Body Text"

# --- Section: Downloads ---------------------------------------------------
$pDownloads = Get-ParaByText "Downloads"
$pDownloads.Range.Text = "##Downloads"
$pDownloads.Style = "Body Text"

$pLeaversChecklist = Get-ParaByText "Leavers checklist"
$pLeaversChecklist.Style = "Body Text"

# --- Section: Contacts -----------------------------------------------------
$pContacts = Get-ParaByText "Contacts"
$pContacts.Range.Text = "##Contacts"
$pContacts.Style = "Body Text"

$pContactsBody = Get-ParaByText "For any further questions relating to security, contact:"
if ($null -eq $pContactsBody) {
    # Text() trims mid-paragraph by exact match above; the paragraph holds
    # several runs, so fall back to a prefix search over all paragraphs.
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith("For any further questions relating to security, contact:")) {
            $pContactsBody = $p
            break
        }
    }
}
$pContactsBody.Style = "Body Text"

# --- Section: Feedback -------------------------------------------------
$pFeedback = Get-ParaByText "Feedback"
$pFeedback.Range.Text = "##Feedback"
$pFeedback.Style = "FirstParagraph"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("If you have any questions or comments about this guidance")) {
        $p.Style = "Body Text"
        break
    }
}
